$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI LR-pair stats for Gnai2-Cnr1 (re-run "following Dr Hou advice"):
# recomputed values for existing sending->target rows, plus two new target
# clusters (sCs) added for each sending cluster (ECs, FAPs, sCs).

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Cnr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 203.7816646666667
$ws.Range("H2").Value = 611.344994
$ws.Range("I2").Value = 0.6667327591988204
$ws.Range("J2").Value = 0.6667327591988205
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.168796
$ws.Range("N2").Value = 6.506387999999999
$ws.Range("O2").Value = 0.9911642957173927
$ws.Range("P2").Value = 0.9911642957173927
$ws.Range("Q2").Value = 441.960859202408
$ws.Range("R2").Value = 3977.647732821672
$ws.Range("S2").Value = 0.6608417057030128
$ws.Range("T2").Value = 0.6608417057030129

# Row 3: ECs -> sCs (new row)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Cnr1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 203.7816646666667
$ws.Range("H3").Value = 611.344994
$ws.Range("I3").Value = 0.6667327591988204
$ws.Range("J3").Value = 0.6667327591988205
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01933366666666667
$ws.Range("N3").Value = 0.058001
$ws.Range("O3").Value = 0.00883570428260726
$ws.Range("P3").Value = 0.008835704282607262
$ws.Range("Q3").Value = 3.939846777443778
$ws.Range("R3").Value = 35.458620996994
$ws.Range("S3").Value = 0.005891053495807572
$ws.Range("T3").Value = 0.005891053495807574

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Cnr1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("H4").Value = 189.421768
$ws.Range("I4").Value = 0.2065833519051582
$ws.Range("J4").Value = 0.2065833519051582
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.168796
$ws.Range("N4").Value = 6.506387999999999
$ws.Range("O4").Value = 0.9911642957173927
$ws.Range("P4").Value = 0.9911642957173927
$ws.Range("Q4").Value = 136.939057583776
$ws.Range("R4").Value = 1232.451518253984
$ws.Range("S4").Value = 0.2047580424980144
$ws.Range("T4").Value = 0.2047580424980144

# Row 5: FAPs -> sCs (new row)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Cnr1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.2065833519051582
$ws.Range("J5").Value = 0.2065833519051582
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01933366666666667
$ws.Range("N5").Value = 0.058001
$ws.Range("O5").Value = 0.00883570428260726
$ws.Range("P5").Value = 0.008835704282607262
$ws.Range("Q5").Value = 1.220739107307556
$ws.Range("R5").Value = 10.986651965768
$ws.Range("S5").Value = 0.001825309407143769
$ws.Range("T5").Value = 0.001825309407143769

# Row 6: sCs -> FAPs (new row)
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Cnr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.719942
$ws.Range("H6").Value = 116.159826
$ws.Range("I6").Value = 0.1266838888960214
$ws.Range("J6").Value = 0.1266838888960214
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.168796
$ws.Range("N6").Value = 6.506387999999999
$ws.Range("O6").Value = 0.9911642957173927
$ws.Range("P6").Value = 0.9911642957173927
$ws.Range("Q6").Value = 83.97565532983199
$ws.Range("R6").Value = 755.7808979684879
$ws.Range("S6").Value = 0.1255645475163655
$ws.Range("T6").Value = 0.1255645475163655

# Row 7: sCs -> sCs (new row)
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Cnr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.719942
$ws.Range("H7").Value = 116.159826
$ws.Range("I7").Value = 0.1266838888960214
$ws.Range("J7").Value = 0.1266838888960214
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01933366666666667
$ws.Range("N7").Value = 0.058001
$ws.Range("O7").Value = 0.00883570428260726
$ws.Range("P7").Value = 0.008835704282607262
$ws.Range("Q7").Value = 0.7485984519806665
$ws.Range("R7").Value = 6.737386067826
$ws.Range("S7").Value = 0.001119341379655919
$ws.Range("T7").Value = 0.001119341379655919
Write-Host "Updated rows 2-7 with new NATMI LR-pair values"
